# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) within specific bullet /
# paragraph text runs, matching the target diff.

$d = $word.ActiveDocument

# RGB 2C3E50 expressed as a Word BGR-ordered color integer (0x50 3E 2C)
$metricColor = 5258796
$pm = [char]0x00B1   # "±"

function Highlight-Metrics {
    param(
        [int]$ParaIndex,
        [string[]]$Targets
    )

    if ($ParaIndex -lt 1) {
        return
    }

    $para = $d.Paragraphs.Item($ParaIndex)
    $searchStart = $para.Range.Start
    $paraEnd = $para.Range.End

    foreach ($target in $Targets) {
        $r = $d.Range($searchStart, $paraEnd)
        $found = $r.Find.Execute($target, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if ($found) {
            $r.Font.Bold = 1
            $r.Font.Color = $metricColor
            $searchStart = $r.End
        }
    }
}

# Locate target paragraphs by matching a distinctive substring, so the
# script is resilient to any paragraph renumbering.
function Find-ParagraphIndex {
    param([string]$Needle)

    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.Contains($Needle)) {
            return $i
        }
    }
    return -1
}

# 1. "Discovered systematic race coding errors ... from 23% to 64%"
$idx = Find-ParagraphIndex "Discovered systematic race coding errors"
Highlight-Metrics $idx @("23%", "64%")

# 2. "Utilized advanced sampling methods ... margin of error from ±4.2% to
#    ±2.1%, ... turnout prediction accuracy from 71% to 87% ..."
$idx = Find-ParagraphIndex "Utilized advanced sampling methods"
Highlight-Metrics $idx @("${pm}4.2%", "${pm}2.1%", "71%", "87%")

# 3. "Trigonometric algorithm ... reduced mapping costs by 73.5%, saving
#    campaigns and organizations $4.7M ..."
$idx = Find-ParagraphIndex "Trigonometric algorithm for boundary estimation"
Highlight-Metrics $idx @("73.5%", "$4.7M")

# 4. "Built real-time FEC analysis systems ... valued over $2 trillion"
$idx = Find-ParagraphIndex "Built real-time FEC analysis systems"
Highlight-Metrics $idx @("$2")

# 5. "Algorithmic innovation: Pioneered trigonometric boundary estimation
#    reducing mapping costs 73.5%"
$idx = Find-ParagraphIndex "Algorithmic innovation: Pioneered trigonometric"
Highlight-Metrics $idx @("73.5%")

# 6. "$4.7M savings enabled nonprofit access"
$idx = Find-ParagraphIndex "savings enabled nonprofit access"
Highlight-Metrics $idx @("$4.7M")

# 7. "178% accuracy improvement in racial classification algorithms"
$idx = Find-ParagraphIndex "accuracy improvement in racial classification"
Highlight-Metrics $idx @("178%")
